$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (LSTM)
$ws.Range("B3").Value = 65.34
$ws.Range("D3").Value = 65.17
$ws.Range("E3").Value = 65.34
$ws.Range("F3").Value = 65.37

# Row 4 (BiLSTM)
$ws.Range("B4").Value = 73.85
$ws.Range("D4").Value = 73.59
$ws.Range("E4").Value = 73.85
$ws.Range("F4").Value = 73.78

# Row 5 (CNN-BiLSTM)
$ws.Range("B5").Value = 78.34
$ws.Range("D5").Value = 77.65
$ws.Range("E5").Value = 78.33
$ws.Range("F5").Value = 78.31

# Row 6 (BERT)
$ws.Range("B6").Value = 80.78
$ws.Range("D6").Value = 81.49
$ws.Range("E6").Value = 80.78
$ws.Range("F6").Value = 80.76

# Row 7 (CompareNet)
$ws.Range("B7").Value = 82.63
$ws.Range("D7").Value = 83.75
$ws.Range("E7").Value = 82.63
$ws.Range("F7").Value = 82.65

# Reset custom row heights on rows 5 and 7 back to default (auto height)
$ws.Rows(5).AutoFit()
$ws.Rows(7).AutoFit()

# Move the active selection
[void]$ws.Range("C11").Select()
